$wb = $excel.ActiveWorkbook

# 1) Prepend "/" to each deep-linking URL in column A (rows 2-19) of the
#    "deep_linking" sheet.
$wsDeep = $wb.Worksheets.Item("deep_linking")
for ($r = 2; $r -le 19; $r++) {
    $cell = $wsDeep.Cells.Item($r, 1)
    $old = $cell.Text
    $cell.Value = "/" + $old
}

# 2) Make "deep_linking" the active/selected sheet, with A4 selected.
$wsDeep.Activate()
$wsDeep.Select()
$wsDeep.Range("A4").Select()
